$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# --- Column width adjustments (column D and column H) ---
# ColumnWidth set via COM is offset by +0.8333... vs the stored XML "width"
# attribute, so subtract that padding to land on the exact target widths.
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668
$ws.Columns.Item(8).ColumnWidth = 12.166666666666666

# --- Row 2 ---
$ws.Range("A2").Value = "2026-01-27 06:32:07"
$ws.Range("B2").Value = "建設・土木業界向け施工機械のAI自動制御・アタッチメント開発を支援してくださるエンジニア募集"
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5434128"
$ws.Range("G2").Value = 368

# --- Row 3 ---
$ws.Range("A3").Value = "2026-01-27 06:32:07"
$ws.Range("B3").Value = "大手製造業向け センサー画像解析・高画質化のR&Dを支援するAIエンジニア募集(画像生成/超解像)"
$ws.Range("D3").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5427956"
$ws.Range("G3").Value = 310
$ws.Range("H3").Value = "🔥AI,Ai"

# --- Row 4 ---
$ws.Range("A4").Value = "2026-01-27 06:32:07"
$ws.Range("B4").Value = "【急募】製造業向け「製造副産物」の状態(硬度)判定AIのフィジビリティ検証(画像認識/動画解析)"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5439158"
$ws.Range("G4").Value = 303
$ws.Range("H4").Value = "🔥AI,Ai"

# --- Row 5 ---
$ws.Range("A5").Value = "2026-01-27 06:32:07"
$ws.Range("B5").Value = "【急募】メルカリ出品商品の在庫管理自動化ツール開発"
$ws.Range("D5").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5479836"
$ws.Range("G5").Value = 215
$ws.Range("H5").Value = "◆ツール,開発 ◇管理"

# --- Row 6 (only the timestamp changes) ---
$ws.Range("A6").Value = "2026-01-27 06:32:07"

# --- Row 7 ---
$ws.Range("A7").Value = "2026-01-27 06:32:07"
$ws.Range("B7").Value = "【北海道・沖縄】2026年度新人Java研修講師募集!3カ月の短期"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5479693"
$ws.Range("G7").Value = 85

# --- Row 8 ---
$ws.Range("A8").Value = "2026-01-27 06:32:07"
$ws.Range("B8").Value = "【急募】パルワールドのMOD開発に関する依頼"
$ws.Range("D8").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5479655"
$ws.Range("G8").Value = 68

# --- Row 9 ---
$ws.Range("A9").Value = "2026-01-27 06:32:07"
$ws.Range("B9").Value = "[日本人限定]webシステム制作のプロジェクト管理業務"
$ws.Range("D9").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5479860"
$ws.Range("G9").Value = 60
$ws.Range("H9").Value = "◇管理"

# --- Row 10 ---
$ws.Range("A10").Value = "2026-01-27 06:32:07"
$ws.Range("B10").Value = "【急募】WordPressでの会議室予約システム構築依頼"
$ws.Range("D10").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5479809"
$ws.Range("G10").Value = 48
$ws.Range("H10").Value = "○WordPress"

# --- Row 11 ---
$ws.Range("A11").Value = "2026-01-27 06:32:07"
$ws.Range("B11").Value = "AntigravityからAndroid、iOSを含めたアプリのリリースを教えてほしい"
$ws.Range("D11").Value = "1,000 ~ 5,000 円 / 固定"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5479715"
$ws.Range("G11").Value = 30
$ws.Range("H11").Value = "◇アプリ"
